$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab ("Лист1" -> "1") ---
$ws.Name = "1"

# --- Fix row 2 data values (validators/goods edits) ---
$ws.Range("A2").Value = "Fdgfbfgb"
$ws.Range("B2").Value = 23344.32
$ws.Range("D2").Value = "dfdbfgb"
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 8
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 4

# --- Add more empty (but formatted) goods rows, 3 through 20 ---
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3:B20").PasteSpecial(-4122) | Out-Null
$ws.Range("K2:O2").Copy() | Out-Null
$ws.Range("K3:O20").PasteSpecial(-4122) | Out-Null
$ws.Range("Q2").Copy() | Out-Null
$ws.Range("Q3:Q20").PasteSpecial(-4122) | Out-Null
$ws.Range("R2").Copy() | Out-Null
$ws.Range("R3:R20").PasteSpecial(-4122) | Out-Null
$ws.Range("T2:U2").Copy() | Out-Null
$ws.Range("T3:U20").PasteSpecial(-4122) | Out-Null
$ws.Range("AM2").Copy() | Out-Null
$ws.Range("AM3:AM20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Fix the "range slider" columns: widen/resize all used columns ---
$ws.Columns.Item(1).ColumnWidth = 5.053385416666667
$ws.Columns.Item(2).ColumnWidth = 19.721354166666668
$ws.Columns.Item(3).ColumnWidth = 8.830729166666666
$ws.Columns.Item(4).ColumnWidth = 9.385416666666666
$ws.Columns.Item(5).ColumnWidth = 9.053385416666666
$ws.Columns.Item(6).ColumnWidth = 34.166666666666664
$ws.Columns.Item(7).ColumnWidth = 11.608072916666666
$ws.Columns.Item(8).ColumnWidth = 9.166666666666666
$ws.Columns.Item(9).ColumnWidth = 5.498697916666667
$ws.Columns.Item(10).ColumnWidth = 11.053385416666666
$ws.Columns.Item(11).ColumnWidth = 12.608072916666666
$ws.Columns.Item(12).ColumnWidth = 11.944010416666666
$ws.Columns.Item(13).ColumnWidth = 5.608072916666667
$ws.Columns.Item(14).ColumnWidth = 5.385416666666667
$ws.Columns.Item(15).ColumnWidth = 5.053385416666667
$ws.Columns.Item(16).ColumnWidth = 14.721354166666666
$ws.Columns.Item(17).ColumnWidth = 10.053385416666666
$ws.Columns.Item(18).ColumnWidth = 14.385416666666666
$ws.Columns.Item(19).ColumnWidth = 13.053385416666666
$ws.Columns.Item(20).ColumnWidth = 14.166666666666666
$ws.Columns.Item(21).ColumnWidth = 15.830729166666666
$ws.Columns.Item(22).ColumnWidth = 14.830729166666666
$ws.Columns.Item(23).ColumnWidth = 11.385416666666666
$ws.Columns.Item(24).ColumnWidth = 12.830729166666666
$ws.Columns.Item(25).ColumnWidth = 12.166666666666666
$ws.Columns.Item(26).ColumnWidth = 10.166666666666666
$ws.Columns.Item(27).ColumnWidth = 11.385416666666666
$ws.Columns.Item(28).ColumnWidth = 18.053385416666668
$ws.Columns.Item(29).ColumnWidth = 15.608072916666666
$ws.Columns.Item(30).ColumnWidth = 7.385416666666667
$ws.Columns.Item(31).ColumnWidth = 17.053385416666668
$ws.Columns.Item(32).ColumnWidth = 13.608072916666666
$ws.Columns.Item(33).ColumnWidth = 9.498697916666666
$ws.Columns.Item(34).ColumnWidth = 6.276041666666667
$ws.Columns.Item(35).ColumnWidth = 8.166666666666666
$ws.Columns.Item(36).ColumnWidth = 16.608072916666668
$ws.Columns.Item(37).ColumnWidth = 19.053385416666668
$ws.Columns.Item(38).ColumnWidth = 18.053385416666668
$ws.Columns.Item(39).ColumnWidth = 6.276041666666667

# --- Fix the view: zoom to 70% and move the selection/active cell ---
$excel.ActiveWindow.Zoom = 70
$ws.Range("A7").Select() | Out-Null

# --- Fix the error handler / page setup (portrait orientation) ---
$ws.PageSetup.Orientation = 1
